# Commit message: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables metadata markers embedded as plain text in the first rows of
# each worksheet used UpperCamelCase attribute names (ObjTablesVersion, Type,
# Id). This updates them to lowerCamelCase (objTablesVersion, type, id) to
# match the new ObjTables convention, without altering anything else about
# the cells (same styles, same cell types - inline/literal strings).

$wb = $excel.ActiveWorkbook

# "!!_Table of contents" sheet: A1 holds the top-level ObjTables version
# marker, A2 holds the sheet/table type marker.
$tocSheet = $wb.Worksheets.Item("!!_Table of contents")
$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents'"

# "!!Model1s" sheet: A1 holds the table type + id marker.
$model1sSheet = $wb.Worksheets.Item("!!Model1s")
$model1sSheet.Range("A1").Value = "!!ObjTables type='Data' id='Model1'"
